$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "Data Sheet 1"
$wb.Worksheets.Item(2).Name = "Data Sheet 2"
$wb.Worksheets.Item(3).Name = "Data Sheet 3"
